$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RebellerNoob / SelectScene" row (row 4) entirely.
$ws.Rows.Item(4).Delete()

# Remove the "CloneScene / Scene2" row (row 2) entirely.
$ws.Rows.Item(2).Delete()

# Update the RelivePos for the villageScene (PioneerNoob) row, now row 2.
$ws.Range("E2").Value = "20,0,60"

# Update the ID for the Demo1 row, now row 3, from 4 to 2.
$ws.Range("B3").Value = 2

# Reflect the final selection used when the data was last edited.
$ws.Range("F5").Select()
